$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Make room for two new rows (6 and 7) by inserting three blank rows above
# the current row 6. This pushes the existing rows 6-13 down to rows 9-16,
# preserving all of their original formatting untouched.
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()

# Row 6: bold header "Notes on the Brazilian adaptation ", highlighted
# yellow across columns A:D.
$ws.Range("A6").Value = "Notes on the Brazilian adaptation "
$ws.Range("A6").Font.Bold = $true
$ws.Range("A6:D6").Interior.Color = 65535

# Row 7: body text "We assumed the same output units as the US model. "
$ws.Range("A7").Value = "We assumed the same output units as the US model. "

# Restore the selection to match the saved workbook state.
$ws.Range("J11").Select()
